# Apply the weekly FlashScore odds refresh:
# - Insert a new row at position 6 (MALAYSIA - SUPER LEAGUE: Perak vs Terengganu),
#   which shifts the former rows 6-8 down to rows 7-9.
# - The former row 6 (MEXICO - LIGA MX) and row 7 (SOUTH KOREA - Gangwon vs
#   Gimcheon Sangmu) land correctly at rows 7 and 8 with no further edits needed.
# - The former row 8 (SOUTH KOREA - Suwon FC vs Seoul) lands at row 9, but its
#   odds columns have been refreshed, so we rewrite that row's values too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Insert()

# Row 6
$ws.Range("A6").Value = "d0ZZLJGq"
$ws.Range("B6").Value = "26/10/2024"
$ws.Range("C6").Value = "06:00"
$ws.Range("D6").Value = "MALAYSIA - SUPER LEAGUE"
$ws.Range("E6").Value = "Perak"
$ws.Range("F6").Value = "Terengganu"
$ws.Range("G6").Value = 3.3
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 1.98
$ws.Range("J6").Value = 3.75
$ws.Range("K6").Value = 2.15
$ws.Range("L6").Value = 2.52
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 11.5
$ws.Range("O6").Value = 1.18
$ws.Range("P6").Value = 3.72
$ws.Range("Q6").Value = 1.65
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.34
$ws.Range("T6").Value = 3.1
$ws.Range("U6").Value = 1.59
$ws.Range("V6").Value = 2.29
$ws.Range("W6").Value = 10
$ws.Range("X6").Value = 16
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 37
$ws.Range("AA6").Value = 22
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 5.9
$ws.Range("AE6").Value = 10
$ws.Range("AF6").Value = 35
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 7.7
$ws.Range("AI6").Value = 9.25
$ws.Range("AJ6").Value = 7.1
$ws.Range("AK6").Value = 15.5
$ws.Range("AL6").Value = 12
$ws.Range("AM6").Value = 17
$ws.Range("AN6").Value = 5.4
$ws.Range("AO6").Value = 18.5
$ws.Range("AP6").Value = 22
$ws.Range("AQ6").Value = 90
$ws.Range("AR6").Value = 110
$ws.Range("AS6").Value = 250
$ws.Range("AT6").Value = 2.95
$ws.Range("AU6").Value = 6.5
$ws.Range("AV6").Value = 50
$ws.Range("AW6").Value = 51
$ws.Range("AX6").Value = 4.05
$ws.Range("AY6").Value = 10
$ws.Range("AZ6").Value = 16
$ws.Range("BA6").Value = 35
$ws.Range("BB6").Value = 55
$ws.Range("BC6").Value = 175
$ws.Range("BD6").Value = 51

# Row 9
$ws.Range("A9").Value = "86Td3Gio"
$ws.Range("B9").Value = "26/10/2024"
$ws.Range("C9").Value = "04:30"
$ws.Range("D9").Value = "SOUTH KOREA - K LEAGUE 1"
$ws.Range("E9").Value = "Suwon FC"
$ws.Range("F9").Value = "Seoul"
$ws.Range("G9").Value = 3.3
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 2.15
$ws.Range("J9").Value = 3.75
$ws.Range("K9").Value = 2.1
$ws.Range("L9").Value = 2.88
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("W9").Value = 10
$ws.Range("X9").Value = 17
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 34
$ws.Range("AA9").Value = 29
$ws.Range("AB9").Value = 34
$ws.Range("AC9").Value = 9
$ws.Range("AD9").Value = 6
$ws.Range("AE9").Value = 15
$ws.Range("AF9").Value = 51
$ws.Range("AG9").Value = 251
$ws.Range("AH9").Value = 7.5
$ws.Range("AI9").Value = 10
$ws.Range("AJ9").Value = 9
$ws.Range("AK9").Value = 21
$ws.Range("AL9").Value = 19
$ws.Range("AM9").Value = 29
$ws.Range("AN9").Value = 5
$ws.Range("AO9").Value = 19
$ws.Range("AP9").Value = 26
$ws.Range("AQ9").Value = 51
$ws.Range("AR9").Value = 81
$ws.Range("AS9").Value = 201
$ws.Range("AT9").Value = 2.63
$ws.Range("AU9").Value = 8
$ws.Range("AV9").Value = 51
$ws.Range("AW9").Value = 501
$ws.Range("AX9").Value = 4.33
$ws.Range("AY9").Value = 12
$ws.Range("AZ9").Value = 23
$ws.Range("BA9").Value = 41
$ws.Range("BB9").Value = 67
$ws.Range("BC9").Value = 151
$ws.Range("BD9").Value = 51

Write-Output "Row insert + odds refresh applied."
